# Automatische test-sync: 2025-07-27 19:27:50
# Adds the new mail-log entry (Testmail #8) as row 10 on the "Logs" sheet,
# extends the conditional-formatting ranges to include the new row, and
# refreshes the "Dashboard" category summary so that "Productinformatie"
# (now 3 occurrences) is listed before "Intern verzoek / Actie voor
# medewerker" (still 2 occurrences).

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(10, 1).Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item(10, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(10, 3).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Cells.Item(10, 4).Value = "Productinformatie"
# Column E (Antwoord) intentionally left blank for this row.
$logs.Cells.Item(10, 6).Value = "2025-07-27 19:27:45"
$logs.Cells.Item(10, 7).Value = "Nee"
$logs.Cells.Item(10, 8).Value = "Ja"
$logs.Cells.Item(10, 9).Value = "Nee"
$logs.Cells.Item(10, 10).Value = "Nee"

# --- 2. Extend the conditional formatting ranges to cover the new row ------
$newRanges = @{
    "D2:D9" = "D2:D10"
    "G2:G9" = "G2:G10"
    "H2:H9" = "H2:H10"
    "I2:I9" = "I2:I10"
    "J2:J9" = "J2:J10"
}

foreach ($oldAddr in $newRanges.Keys) {
    $newAddr = $newRanges[$oldAddr]
    $rng = $logs.Range($oldAddr)
    for ($i = 1; $i -le $rng.FormatConditions.Count(); $i++) {
        $fc = $rng.FormatConditions.Item($i)
        $fc.ModifyAppliesToRange($logs.Range($newAddr))
    }
}

# --- 3. Update the "Dashboard" category summary table -----------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(3, 1).Value = "Productinformatie"
$dashboard.Cells.Item(3, 2).Value = 3
$dashboard.Cells.Item(4, 1).Value = "Intern verzoek / Actie voor medewerker"
$dashboard.Cells.Item(4, 2).Value = 2
